$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rota")

# New "Own Jobs" worktype roster data added to the export.
# Column A = OFF, B = ON, C = HOUSECARE, D = DEBOP, F = OWN JOBS, H = KITCHEN

$ws.Range("A3").Value = "Adam PM"
$ws.Range("B3").Value = "Navi"
$ws.Range("C3").Value = "MahaDeva"
$ws.Range("F3").Value = "Mahi - OJ"
$ws.Range("H3").Value = "Dganit"

$ws.Range("A4").Value = "Supriti"
$ws.Range("B4").Value = "Neal"
$ws.Range("F4").Value = "Khalsa - PR"
$ws.Range("H4").Value = "Ben"

$ws.Range("A5").Value = "Anna AM"
$ws.Range("B5").Value = "Mel"
$ws.Range("F5").Value = "Adam - BY"

$ws.Range("A6").Value = "Shakti"
$ws.Range("B6").Value = "MahaDeva"

$ws.Range("A7").Value = "Anuka AM"
$ws.Range("B7").Value = "Khalsa"

$ws.Range("B8").Value = "Adam AM"
$ws.Range("B9").Value = "Anna PM"
$ws.Range("B10").Value = "Dganit"
$ws.Range("B11").Value = "Anuka PM"
$ws.Range("B12").Value = "Ben"
$ws.Range("B13").Value = "Mahi"
$ws.Range("B14").Value = "RAP"

$ws.Range("C18").Value = ""

$ws.Range("C19").Value = "Dganit"
$ws.Range("D19").Value = "MahaDeva"
$ws.Range("F19").Value = "Anuka - OJ"
$ws.Range("H19").Value = "Mahi"

$ws.Range("C20").Value = "Anna"
$ws.Range("F20").Value = "Khalsa - PR"
$ws.Range("H20").Value = "Ben"

$ws.Columns.Item(1).ColumnWidth = 9.692829767862966
$ws.Columns.Item(4).ColumnWidth = 10.035557428995766

